$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Valid sets")

# Widen column A and drop the bestFit auto-sizing (now a fixed custom width)
$ws.Columns.Item(1).ColumnWidth = 41.333333

# Insert two new rows before row 22 (pushing the old row22/23 down to 24/25)
$ws.Rows.Item(22).Resize(2).Insert() | Out-Null

# Fill new row 23 first (3-line string) so it becomes shared string index 29
$ws.Range("A23").Value = "training-run3-test75-20210924-0517`ntraining-run3-test75-20210924-0717`ntraining-run3-test75-20210924-0817"
$ws.Range("B23").Value = "V6"
$ws.Range("C23").Value = 903632
$ws.Range("D23").Value = 98271032

# Fill new row 22 (2-line string) so it becomes shared string index 30
$ws.Range("A22").Value = "training-run3-test75-20210924-0117`ntraining-run3-test75-20210924-0217"
$ws.Range("B22").Value = "V6"
$ws.Range("C22").Value = 600000
$ws.Range("D22").Value = 65242480

# Wrap the text in the two new multi-line name cells, and size their row height
# to match the 2-line / 3-line wrapped content
$ws.Range("A22:A23").WrapText = $true
$ws.Rows.Item(22).RowHeight = 29
$ws.Rows.Item(23).RowHeight = 43.5

# Recompute E formulas for every data row (now 2..25), cumulative sum
$ws.Range("E2").Formula = "=D2"
$ws.Range("E3:E25").Formula = "=E2+D3"

# New column F: "Total number of chunks" cumulative sum of C
$ws.Range("F1").Value = "Total number of chunks"
$ws.Range("F2").Formula = "=C2"
$ws.Range("F3:F25").Formula = "=F2+C3"

# Update the averaged-positions-per-chunk formula to cover the extended range
$ws.Range("G2").Formula = "=SUM(D2:D25)/SUM(C2:C25)"

# Restore the selection to what it was left at in the authored workbook
$ws.Range("G29").Select() | Out-Null

Write-Host "done"
